$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.566.00'
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").Value = '1.763.86'
$ws.Range("E3").Value = '  -1.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3839'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3412'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.94'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.137'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07395'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.43'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.341'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.72%  '
$ws.Range("D15").Value = '1.764.16'
$ws.Range("E15").Value = '  -1.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.067'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001071'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06667'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.97%  '
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.380'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.37%  '
$ws.Range("D23").Value = '27.565.15'
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.386'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.427'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.60'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.415'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '152.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("D31").Value = '1.963.85'
$ws.Range("E31").Value = '  -1.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.098'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.964'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.28%  '
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.72'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02408'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6779'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.305'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06317'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2176'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.510'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.242'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.249'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.12%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.12'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6238'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.825'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.56'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.070'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07378'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.145'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.49%  '
